$p = $ppt.ActivePresentation
$s = $p.Slides.Item(22)
$sh = $s.Shapes.Item(1)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# "Елементите в дървото са " -> "Елементите са "
$tr.Characters(91, 24).Text = "Елементите са "

# " на дървото връща елементите във " -> " връща елементите във "
$tr.Characters(134, 33).Text = " връща елементите във "
